# Apply updated per-player statistic values (cfs_6_0.4.xlsx data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Iga Swiatek
$ws.Range("H2").Value = 0.407
$ws.Range("I2").Value = 0.639

# Row 3: Aryna Sabalenka
$ws.Range("J3").Value = 0.368

# Row 4: Coco Gauff
$ws.Range("G4").Value = 0.734

# Row 5: Elena Rybakina
$ws.Range("G5").Value = 0.789

# Row 6: Jessica Pegula
$ws.Range("G6").Value = 0.723
$ws.Range("K6").Value = 0.496

# Row 7: Ons Jabeur
$ws.Range("H7").Value = 0.4

# Row 8: Marketa Vondrousova
$ws.Range("G8").Value = 0.707

# Row 9: Karolina Muchova
$ws.Range("C9").Value = 0.663
$ws.Range("L9").Value = 0.446

# Row 10: Maria Sakkari
$ws.Range("H10").Value = 0.367

# Row 11: Barbora Krejcikova
$ws.Range("E11").Value = 0.601

# Row 13: Madison Keys
$ws.Range("L13").Value = 0.447

# Row 14: Jelena Ostapenko
$ws.Range("E14").Value = 0.578
$ws.Range("G14").Value = 0.6860000000000001

# Row 15: Petra Kvitova
$ws.Range("I15").Value = 0.547
$ws.Range("L15").Value = 0.424

# Row 16: Qinwen Zheng
$ws.Range("G16").Value = 0.742

# Row 17: Liudmila Samsonova
$ws.Range("L17").Value = 0.428

# Row 18: Belinda Bencic
$ws.Range("B18").Value = 0.605

# Row 21: Caroline Garcia
$ws.Range("D21").Value = 0.495
$ws.Range("G21").Value = 0.8080000000000001

# Row 23: Victoria Azarenka
$ws.Range("B23").Value = 0.638
$ws.Range("K23").Value = 0.454

# Row 30: Elise Mertens
$ws.Range("H30").Value = 0.391

# Row 31: Jasmine Paolini
$ws.Range("H31").Value = 0.47

# Row 32: Lesia Tsurenko
$ws.Range("E32").Value = 0.524
$ws.Range("I32").Value = 0.574

# Row 43: Martina Trevisan
$ws.Range("K43").Value = 0.359

# Row 44: Varvara Gracheva
$ws.Range("C44").Value = 0.633
$ws.Range("I44").Value = 0.547

# Row 46: Katerina Siniakova
$ws.Range("G46").Value = 0.638

# Row 48: Sloane Stephens
$ws.Range("G48").Value = 0.645
$ws.Range("K48").Value = 0.456

# Row 49: Sara Sorribes Tormo
$ws.Range("C49").Value = 0.588

# Row 56: Clara Burel
$ws.Range("E56").Value = 0.547

# Row 58: Mirra Andreeva
$ws.Range("G58").Value = 0.735

# Row 60: Mayar Sherif
$ws.Range("E60").Value = 0.545

# Row 61: Greet Minnen
$ws.Range("J61").Value = 0.341

# Row 62: Cristina Bucsa
$ws.Range("G62").Value = 0.615

# Row 65: Paula Badosa
$ws.Range("C65").Value = 0.6879999999999999
$ws.Range("E65").Value = 0.603

# Row 68: Nadia Podoroska
$ws.Range("H68").Value = 0.424
$ws.Range("I68").Value = 0.591

# Row 69: Lauren Davis
$ws.Range("G69").Value = 0.645
$ws.Range("K69").Value = 0.456

# Row 71: Xiyu Wang
$ws.Range("K71").Value = 0.336

# Row 73: Anna Karolina Schmiedlova
$ws.Range("E73").Value = 0.511

# Row 76: Irina-Camelia Begu
$ws.Range("D76").Value = 0.461

# Row 77: Anna Kalinskaya
$ws.Range("G77").Value = 0.712

# Row 79: Camila Osorio
$ws.Range("F79").Value = 0.509
$ws.Range("L79").Value = 0.43

# Row 80: Diane Parry
$ws.Range("I80").Value = 0.5629999999999999

# Row 81: Viktoriya Tomova
$ws.Range("G81").Value = 0.616
$ws.Range("I81").Value = 0.537

# Row 82: Taylor Townsend
$ws.Range("G82").Value = 0.772

# Row 84: Viktorija Golubic
$ws.Range("J84").Value = 0.412

# Row 86: Alycia Parks
$ws.Range("I86").Value = 0.547

# Row 87: Laura Siegemund
$ws.Range("K87").Value = 0.337
$ws.Range("L87").Value = 0.387

# Row 88: Kayla Day
$ws.Range("G88").Value = 0.6830000000000001

# Row 89: Zhu Oxuanbai
$ws.Range("I89").Value = 0.609

# Row 90: Jaqueline Cristian
$ws.Range("C90").Value = 0.591
$ws.Range("K90").Value = 0.381

# Row 91: Aliaksandra Sasnovich
$ws.Range("C91").Value = 0.571
$ws.Range("E91").Value = 0.515

# Row 93: Linda Fruhvirtova
$ws.Range("J93").Value = 0.474

# Row 95: Nao Hibino
$ws.Range("G95").Value = 0.671

# Row 97: Oceane Dodin
$ws.Range("C97").Value = 0.641
$ws.Range("G97").Value = 0.647

# Row 98: Jodie Burrage
$ws.Range("H98").Value = 0.395

# Row 99: Claire Liu
$ws.Range("D99").Value = 0.483
$ws.Range("M99").Value = 5.2

# Row 100: Kamilla Rakhimova
$ws.Range("G100").Value = 0.671
